# Generate Report for handoff
# This script updates the localization-status workbook to reflect a failed
# handoff transform for file 17a7fd45-f2b8-4a1e-bf91-9439f2fb2369.md, which
# is replaced by a new report e763d729-3c0a-492f-a7c8-9d63441ec1c3.md, with
# status "Handoff transform failed" and the per-language handoff details
# reset/ignored.

$wb = $excel.ActiveWorkbook

$oldFile = "17a7fd45-f2b8-4a1e-bf91-9439f2fb2369.md"
$newFile = "e763d729-3c0a-492f-a7c8-9d63441ec1c3.md"
$oldStatus = "Ready for handoff"
$newStatus = "Handoff transform failed"
$zeroDate = "0001-01-01 00:00:00"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# --- Overview sheet ---
$wsOverview.Range("A2").Value = $newFile
$wsOverview.Hyperlinks.Item(1).TextToDisplay = $newFile
$wsOverview.Range("B2").Value = $newStatus
$wsOverview.Range("C2").Value = $newStatus

# --- zh-cn sheet ---
$wsZhCn.Range("A2").Value = $newFile
$wsZhCn.Hyperlinks.Item(1).TextToDisplay = $newFile
$wsZhCn.Range("B2").Value = $newStatus
$wsZhCn.Range("C2").Hyperlinks.Delete()
$wsZhCn.Range("C2").Value = ""
$wsZhCn.Range("D2").Value = $zeroDate
$wsZhCn.Range("H2").Value = "Ignored"

# --- de-de sheet ---
$wsDeDe.Range("A2").Value = $newFile
$wsDeDe.Hyperlinks.Item(1).TextToDisplay = $newFile
$wsDeDe.Range("B2").Value = $newStatus
$wsDeDe.Range("C2").Hyperlinks.Delete()
$wsDeDe.Range("C2").Value = ""
$wsDeDe.Range("D2").Value = $zeroDate
$wsDeDe.Range("H2").Value = "Ignored"
